# cn-#14 add ignoreTableFormat option when comparing files.
# Update the "SumProduct" worksheet (sheet2) of the workbook:
#  - Add a third list of numbers (column D) and a per-row product (column E)
#  - Update the SUMPRODUCT formula in C11 to also take column D into account
#  - Update selection/dimension bookkeeping to match

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SumProduct")

# Row 5: third factor + product formula
$ws.Range("D5").Value = 2
$ws.Range("E5").Formula = "=B5*C5*D5"

# Row 6: third factor + shared product formula
$ws.Range("D6").Value = 12
$ws.Range("E6").Formula = "=B6*C6*D6"

# Row 7: drop the old C7 value, add third factor + shared product formula
$ws.Range("C7").ClearContents()
$ws.Range("D7").Value = 1
$ws.Range("E7").Formula = "=B7*C7*D7"

# Update the total formula to include the new column D range
$ws.Range("C11").Formula = "=SUMPRODUCT(B5:B7,C5:C7, D5:D7)"

# Match the final selection recorded in the saved file
$ws.Range("C7").Select()
